# Apply "Better validations for columns and blank rows" edit to the
# Instructions sheet of the workbook.
#
# Summary of the change:
#   - Bump the version string in A2.
#   - Split the old A4 instruction text into three separate lines (A4..A6).
#   - Insert two new rows so the rest of the content (the field reference
#     table starting at what used to be row 6) shifts down by two rows.

$wb = $excel.ActiveWorkbook

# The edited sheet is "Instructions" (the first sheet / sheetId 1).
$ws = $wb.Worksheets.Item("Instructions")

# The sheet is protected; unprotect it so we can edit cell contents and
# insert rows, matching what Excel would require interactively.
$ws.Unprotect()

# Bump version number.
$ws.Range("A2").Value = "Version 1.2.2"

# Insert two blank rows above the old row 6 ("Antibody name" header row),
# pushing everything from there on down by two rows.
$ws.Rows("5:6").Insert()

# Re-split the old combined instructions line (formerly A4) into three
# separate lines.
$ws.Range("A4").Value = "Add your antibodies to the 'Antibodies' sheet."
$ws.Range("A5").Value = "Do not change the headers of the 'Antibodies' sheet."
$ws.Range("A6").Value = "Do not edit the other sheets."

# Restore sheet protection (with default options, no password) to match
# the original protected state.
$ws.Protect()
